$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 267119.850555793
$ws.Range("D2").Value = 287634.092223972
$ws.Range("E2").Value = 280533.392199001
$ws.Range("F2").Value = 253706.308912585
$ws.Range("G2").Value = 246605.608887614
$ws.Range("C3").Value = 254065.399322366
$ws.Range("D3").Value = 277657.796861272
$ws.Range("E3").Value = 269491.638906065
$ws.Range("F3").Value = 238639.159738667
$ws.Range("G3").Value = 230473.001783459
$ws.Range("C4").Value = 292341.960478844
$ws.Range("D4").Value = 323486.745533066
$ws.Range("E4").Value = 312706.440826249
$ws.Range("F4").Value = 271977.480131438
$ws.Range("G4").Value = 261197.175424621
$ws.Range("C5").Value = 282403.506319702
$ws.Range("D5").Value = 315911.660438468
$ws.Range("E5").Value = 304313.310636587
$ws.Range("F5").Value = 260493.702002818
$ws.Range("G5").Value = 248895.352200937
$ws.Range("C6").Value = 318611.407706171
$ws.Range("D6").Value = 359919.739120796
$ws.Range("E6").Value = 345621.47381556
$ws.Range("F6").Value = 291601.341596782
$ws.Range("G6").Value = 277303.076291545
$ws.Range("C7").Value = 326052.780579756
$ws.Range("D7").Value = 371632.725656654
$ws.Range("E7").Value = 355855.904734837
$ws.Range("F7").Value = 296249.656424674
$ws.Range("G7").Value = 280472.835502858
$ws.Range("C8").Value = 329543.068974335
$ws.Range("D8").Value = 378728.702111271
$ws.Range("E8").Value = 361703.825828107
$ws.Range("F8").Value = 297382.312120564
$ws.Range("G8").Value = 280357.4358374
$ws.Range("C9").Value = 295579.204364193
$ws.Range("D9").Value = 342326.853808094
$ws.Range("E9").Value = 326145.849365285
$ws.Range("F9").Value = 265012.559363101
$ws.Range("G9").Value = 248831.554920292
$ws.Range("C10").Value = 265215.831658795
$ws.Range("D10").Value = 309397.847766281
$ws.Range("E10").Value = 294104.899191713
$ws.Range("F10").Value = 236326.764125877
$ws.Range("G10").Value = 221033.81555131
$ws.Range("C11").Value = 260183.194339197
$ws.Range("D11").Value = 305616.516123288
$ws.Range("E11").Value = 289890.446676263
$ws.Range("F11").Value = 230475.942002131
$ws.Range("G11").Value = 214749.872555106
$ws.Range("C12").Value = 268877.697933027
$ws.Range("D12").Value = 317895.15720478
$ws.Range("E12").Value = 300928.491806913
$ws.Range("F12").Value = 236826.904059141
$ws.Range("G12").Value = 219860.238661274
$ws.Range("C13").Value = 260865.151021829
$ws.Range("D13").Value = 310346.437192453
$ws.Range("E13").Value = 293219.225003451
$ws.Range("F13").Value = 228511.077040206
$ws.Range("G13").Value = 211383.864851204
